$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two match records (rows 86 and 87, columns B through AC) had their
# data swapped, while the leading index column A stayed in place.
$range86 = $ws.Range("B86:AC86")
$range87 = $ws.Range("B87:AC87")

$values86 = $range86.Value()
$values87 = $range87.Value()

$range86.Value = $values87
$range87.Value = $values86
